# previsao_retorno.xlsx update ("atualizei dados bibi e add")
#
# Two independent kinds of changes land on the single sheet
# (Resumo_por_Cliente):
#
#  1) Every "INATIVO - X.Y meses sem comprar" status label (column J,
#     "situacao") is the result of a recompute against a later reference
#     date, so each already-inactive client's elapsed-months figure moved
#     forward by exactly 0.1 months.  "ATIVO" rows are untouched.
#
#  2) A handful of specific clients (rows 86, 91, 113, 197, 234, 249, 359,
#     377, 392, 414) got fresh purchase activity, shifting their
#     probabilities (B/C/D/F), purchase counts (E), pattern label (G) and
#     last/next purchase dates (H/I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bump every "INATIVO - X.Y meses sem comprar" label by +0.1 months.
# ---------------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)   # column J = situacao
    $text = $cell.Text
    if ($text -match '^INATIVO - ([\d\.]+) meses sem comprar$') {
        $meses = [double]$matches[1] + 0.1
        $cell.Value2 = "INATIVO - {0:N1} meses sem comprar" -f $meses
    }
}

# ---------------------------------------------------------------------------
# 2) Per-client refreshed figures.
# ---------------------------------------------------------------------------
$rowUpdates = @(
    @{ Row = 86;  B = 0.83; C = 0.83;              E = 33;  H = 45856; I = 45918 },
    @{ Row = 91;                                   E = 312; H = 45859; I = 45866 },
    @{ Row = 113; B = 0.58; C = 0.17;              E = 99;  H = 45856; I = 45871 },
    @{ Row = 197;                                  E = 52;  H = 45856; I = 45887; G = '1x por mês - irregular (preferencialmente na 2ª quinzena)' },
    @{ Row = 234; B = 0.58; C = 0.33;              E = 32;  H = 45856; I = 45918 },
    @{ Row = 249; B = 0.75; C = 0.5;                         H = 45856; I = 45871 },
    @{ Row = 359; B = 0.58;           D = 0.67; F = 0.67; E = 18;  H = 45856; I = 45887 },
    @{ Row = 377; B = 0.08;           D = 0.17; F = 0.17; E = 9;   H = 45845; I = 45907; G = '1x a cada 2 meses - irregular' },
    @{ Row = 392; B = 0.75;           D = 1;    F = 1;    E = 35;  H = 45856; I = 45871 },
    @{ Row = 414; B = 0.75; C = 0.67;                       H = 45853; I = 45860 }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value2 = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value2 = $u.C }
    if ($u.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value2 = $u.D }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value2 = $u.E }
    if ($u.ContainsKey('F')) { $ws.Cells.Item($r, 6).Value2 = $u.F }
    if ($u.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value2 = $u.G }
    if ($u.ContainsKey('H')) { $ws.Cells.Item($r, 8).Value2 = $u.H }
    if ($u.ContainsKey('I')) { $ws.Cells.Item($r, 9).Value2 = $u.I }
}

Write-Host "Done."
